# Commit: "Doing Updates for Financials"
# Insert two new quarterly columns (D, E) into the AIRT sheet to hold the
# latest two quarters of financial data. Excel's column insert shifts the
# existing D:K data right to F:M, carrying along its formatting; we then
# populate the two new columns with the newly reported figures and restore
# the same per-row number formats (date header row vs. general numeric
# rows) that the rest of the row already uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AIRT")

# Insert two blank columns before column D (old D:K -> F:M)
$ws.Range("D:E").Insert()

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 63600
$ws.Range("E8").Value = 49300
$ws.Range("D9").Value = 49400
$ws.Range("E9").Value = 40000
$ws.Range("D10").Value = 14200
$ws.Range("E10").Value = 9300
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 2000
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 2300
$ws.Range("E15").Value = 1800
$ws.Range("D17").Value = 64200
$ws.Range("E17").Value = 51000
$ws.Range("D18").Value = -600
$ws.Range("E18").Value = -1700
$ws.Range("D20").Value = -1600
$ws.Range("E20").Value = -200
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = 0
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("D23").Value = -2100
$ws.Range("E23").Value = -1800
$ws.Range("D24").Value = 200
$ws.Range("E24").Value = -400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -2300
$ws.Range("E26").Value = -1400
$ws.Range("D27").Value = -2700
$ws.Range("E27").Value = -1300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 1600
$ws.Range("E32").Value = 200
$ws.Range("D33").Value = -2700
$ws.Range("E33").Value = -1300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -2700
$ws.Range("E35").Value = -1300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 2700
$ws.Range("E41").Value = 5600
$ws.Range("D42").Value = 2500
$ws.Range("E42").Value = 2300
$ws.Range("D43").Value = 24100
$ws.Range("E43").Value = 24100
$ws.Range("D44").Value = 39600
$ws.Range("E44").Value = 29500
$ws.Range("D45").Value = 2700
$ws.Range("E45").Value = 3000
$ws.Range("D46").Value = 71600
$ws.Range("E46").Value = 64500
$ws.Range("D47").Value = 6600
$ws.Range("E47").Value = 9200
$ws.Range("D48").Value = 31800
$ws.Range("E48").Value = 37100
$ws.Range("D49").Value = 5700
$ws.Range("E49").Value = 5800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1200
$ws.Range("E52").Value = 1200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 116900
$ws.Range("E54").Value = 117700
$ws.Range("D57").Value = 11000
$ws.Range("E57").Value = 14700
$ws.Range("D58").Value = 33400
$ws.Range("E58").Value = 15200
$ws.Range("D59").Value = 12200
$ws.Range("E59").Value = 13800
$ws.Range("D60").Value = 56600
$ws.Range("E60").Value = 43800
$ws.Range("D61").Value = 33400
$ws.Range("E61").Value = 44100
$ws.Range("D62").Value = 1700
$ws.Range("E62").Value = 1500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 93600
$ws.Range("E66").Value = 91000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 18700
$ws.Range("E72").Value = 22100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 23300
$ws.Range("E76").Value = 26700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -2700
$ws.Range("E81").Value = -1300
$ws.Range("D83").Value = 2200
$ws.Range("E83").Value = 1800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -13600
$ws.Range("E89").Value = 200
$ws.Range("D91").Value = 19000
$ws.Range("E91").Value = -19500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 4800
$ws.Range("E94").Value = -20900
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 6700
$ws.Range("E100").Value = 20600
$ws.Range("D101").Value = 100
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -2000
$ws.Range("E102").Value = -100

# Restore formatting for the newly inserted columns: the "Period Ending"
# rows use the date style, everything else uses the general-number style
# that already governs the rest of each row.
$dateRows = @(7, 38, 80)
foreach ($r in $dateRows) {
    $rng = $ws.Range("D$r" + ":E$r")
    $rng.NumberFormat = "[$-409]d\-mmm\-yy;@"
    $rng.Font.Name = "Verdana"
    $rng.Font.Size = 12
    $rng.Font.Bold = $true
}

$numberRows = @(8..35) + @(39..77) + @(81..102)
foreach ($r in $numberRows) {
    $rng = $ws.Range("D$r" + ":E$r")
    $rng.NumberFormat = "#,##0"
    $rng.Font.Name = "Verdana"
    $rng.Font.Size = 12
    $rng.Font.Bold = $false
    $rng.HorizontalAlignment = -4152
}
